# "Adding basic extent report"
#
# The "login" test-data sheet is replaced with a smaller, refreshed sheet
# (header row + a single "random1" row + the "testemail" row) and becomes
# the active/selected tab. The "register" sheet is left as-is (still the
# duplicate testemail/testemail rows) but is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$registerSheet = $wb.Worksheets.Item("register")
$oldLoginSheet = $wb.Worksheets.Item("login")

# Drop the old 6-row login sheet entirely ...
$oldLoginSheet.Delete()

# ... and rebuild it fresh, positioned right after "register".
$loginSheet = $wb.Worksheets.Add($null, $registerSheet)
$loginSheet.Name = "login"

$loginSheet.Range("A1:B2").ColumnWidth = 18.77734375

# Header row
$loginSheet.Cells.Item(1, 1).Value = "email"
$loginSheet.Cells.Item(1, 2).Value = "password"

# Row 2: random1
$loginSheet.Cells.Item(2, 2).Value = "random1"
$h2 = $loginSheet.Hyperlinks.Add($loginSheet.Range("A2"), "mailto:random1@gmail.com", $null, $null, "random1@gmail.com")

# Row 3: testemail
$loginSheet.Cells.Item(3, 2).Value = "testemail"
$h3 = $loginSheet.Hyperlinks.Add($loginSheet.Range("A3"), "mailto:testemail@gmail.com", $null, $null, "testemail@gmail.com")

# "login" becomes the active tab, with B2 selected.
$loginSheet.Activate()
$loginSheet.Range("B2").Select()
